# Issue #70: Add the calculation of the confidence intervals to the files
# with results (#76)
#
# Adds three new columns (E: Estimate_on_the_response_scale,
# F: lwr_on_the_response_scale, G: uppr_on_the_response_scale) to both
# worksheets ("POP" and "Pesticide"). These hold the model estimate and
# confidence-interval bounds back-transformed onto the response scale
# (i.e. exp() of the Value/Std.-Error columns), rounded to 4 decimals,
# exactly as produced upstream. A couple of the back-transformed upper
# bounds overflow double precision and are recorded as the literal
# #NUM! error, matching what Excel itself would compute.

$wb = $excel.ActiveWorkbook

$headers = @(
    "Estimate_on_the_response_scale",
    "lwr_on_the_response_scale",
    "uppr_on_the_response_scale"
)

# Row layout: @(row, Estimate_on_the_response_scale, lwr_on_the_response_scale, uppr_on_the_response_scale)
$popRows = @(
    @(2, 8.4387, 4.6088, 15.4514),
    @(3, 1.4208, 1.0076, 2.0036),
    @(4, 0.8401, 0.467, 1.5111),
    @(5, 1.3319, 0.9331, 1.901),
    @(6, 0.7298, 0.2379, 2.2382),
    @(7, 1.4626, 1.0052, 2.1282),
    @(8, 1.1565, 0.8014, 1.6689),
    @(9, 0.7422, 0.2302, 2.393),
    @(10, 1.4247, 0.73, 2.7807),
    @(11, 0.9548, 0.6116, 1.4905),
    @(12, 0.9359, 0.5297, 1.6538),
    @(13, 1.0049, 0.5792, 1.7435),
    @(14, 1.0034, 0.5991, 1.6805),
    @(15, 0.9093, 0.5439, 1.5204),
    @(16, 0.8719, 0.5015, 1.5161),
    @(17, 0.8332, 0.4395, 1.5797),
    @(18, 0.7653, 0.3642, 1.6083),
    @(19, 0.7431, 0.3329, 1.6586),
    @(20, 0.7827, 0.3355, 1.8263),
    @(21, 0.7981, 0.3051, 2.0876),
    @(22, 0.808, 0.2277, 2.8668),
    @(23, 0.267, 0.2279, 0.3129)
)

$pesticideRows = @(
    @(2, 0.3499, 0, 232173.9463),
    @(3, 0.0667, 0, "#NUM!"),
    @(4, 1.6801, 0.8581, 3.2893),
    @(5, 0.0667, 0, "#NUM!"),
    @(6, 1.5169, 0.8051, 2.8578),
    @(7, 0.2099, 0, "#NUM!"),
    @(8, 0.0481, 0, "#NUM!"),
    @(9, 0.0548, 0, "#NUM!"),
    @(10, 1.0492, 0.1524, 7.2219),
    @(11, 1.5513, 0.0034, 709.1733),
    @(12, 2.4067, 0.0001, 91787.2196),
    @(13, 3.7251, 0, 1725734.2226),
    @(14, 5.8266, 0, 4556282.3575),
    @(15, 7.6143, 0, 5130972.8681),
    @(16, 9.0815, 0, 6247668.1047),
    @(17, 10.9358, 0, 7914255.3336),
    @(18, 15.7876, 0, 12510838.8833),
    @(19, 8.2196, 0, 6653993.4497),
    @(20, 2.4679, 0, 2757884.0083),
    @(21, 0.7369, 0, 3149435.6347),
    @(22, 0.22, 0, 16344574.4757),
    @(23, 0.5427, 0.3818, 0.7715)
)

function Add-ConfidenceIntervalColumns {
    param(
        $Worksheet,
        $Rows
    )

    # Headers on row 1, columns E (5), F (6), G (7)
    for ($col = 0; $col -lt $headers.Length; $col++) {
        $Worksheet.Cells.Item(1, 5 + $col).Value = $headers[$col]
    }

    foreach ($entry in $Rows) {
        $rowNum = $entry[0]
        for ($col = 0; $col -lt 3; $col++) {
            $Worksheet.Cells.Item($rowNum, 5 + $col).Value = $entry[1 + $col]
        }
    }
}

$wsPop = $wb.Worksheets.Item("POP")
Add-ConfidenceIntervalColumns $wsPop $popRows

$wsPesticide = $wb.Worksheets.Item("Pesticide")
Add-ConfidenceIntervalColumns $wsPesticide $pesticideRows
